# daily auto push: 2026-02-15 03:21 UTC
# Insert a new daily data point for 2026/02/15 (日) just before the
# 2026/12/29 block, shifting all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 823..864 down to 824..865 by inserting a fresh row at 823.
$ws.Rows.Item(823).Insert()

# Column A holds plain text dates (e.g. "2026/02/15"), not real Excel
# dates. Force text format first so the string literal isn't silently
# reinterpreted as a date serial number.
$ws.Cells.Item(823, 1).NumberFormat = "@"
$ws.Cells.Item(823, 1).Value2 = "2026/02/15"
$ws.Cells.Item(823, 2).Value2 = "日"
$ws.Cells.Item(823, 3).Value2 = 8
$ws.Cells.Item(823, 4).Value2 = 201
